$wb = $excel.ActiveWorkbook

$wsSettings  = $wb.Worksheets.Item("Settings")
$wsConstants = $wb.Worksheets.Item("Constants")

# ---------------------------------------------------------------------------
# Constants sheet: move "downloadPath" down, and fold the old "folder does
# not exist" message into the new "logMessage*" naming scheme, then add the
# new cv/job download paths, log messages and reply messages.
# ---------------------------------------------------------------------------
$wsConstants.Range("A9").Value = "downloadPath"
$wsConstants.Range("B9").Value = "Data\Downloads"

$wsConstants.Range("B2").Value = "The folder does not exist. It was created automatically by the system."

$wsConstants.Range("A10").Value = "cvDownloadPath"
$wsConstants.Range("A11").Value = "jobDownloadPath"
$wsConstants.Range("B11").Value = "Data\Downloads\Jobs"
$wsConstants.Range("B10").Value = "Data\Downloads\CVs"

# ---------------------------------------------------------------------------
# Settings sheet: replace the old "downloadPath" row with "email" / "emailFolder"
# ---------------------------------------------------------------------------
$wsSettings.Range("A2").Value = "email"
$wsSettings.Range("B2").Value = "rpa_project@outlook.com"
$wsSettings.Range("A3").Value = "emailFolder"
$wsSettings.Range("B3").Value = "Inbox"

# Turn the email address into a real mailto: hyperlink (gives it the
# built-in "Hyperlink" style too).
$wsSettings.Hyperlinks.Add($wsSettings.Range("B2"), "mailto:rpa_project@outlook.com") | Out-Null

$wsSettings.Columns.Item(1).ColumnWidth = 15.6667
$wsSettings.Columns.Item(2).ColumnWidth = 23.0

$wsSettings.Range("A4").Select()

# ---------------------------------------------------------------------------
# Constants sheet: finish filling in the remaining new rows.
# ---------------------------------------------------------------------------
$wsConstants.Range("A2").Value = "logMessageNewFolder"

$wsConstants.Range("A3").Value = "logMessageNewFile"
$wsConstants.Range("B3").Value = "The file does not exist. It was created automatically by the system."

$wsConstants.Range("A6").Value = "logMessageReadEmails"
$wsConstants.Range("B6").Value = "Reading the emails…"

$wsConstants.Range("A7").Value = "logMessageVerifyExistingFolder"
$wsConstants.Range("B7").Value = "Verifying existence of folder…"

$wsConstants.Range("A5").Value = "logMessageInitAllSettings"
$wsConstants.Range("B5").Value = "Initializing all settings…"

$wsConstants.Range("A13").Value = "replyMessageNoAttachment"
$wsConstants.Range("A14").Value = "replyMessageConfirmation"
$wsConstants.Range("B14").Value = "We have received the email sent by you."
$wsConstants.Range("B13").Value = "The mail you sent does not include an attachment."

$wsConstants.Columns.Item(1).ColumnWidth = 28.5

$wsConstants.Range("D19").Select()
